# Commit: "Sat, May 23, 2020 10:05:01 AM"
#
# The underlying change swaps the deck's two theme parts (the active
# "Integral"/Red Violet theme <-> the spare "Office Theme") and restyles
# the three tables (slides 14-16) from the custom "Table_0" style to the
# built-in table style {C5A9637D-327B-4B98-B6D2-77B0DE613E2E}.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables (slides 14, 15, 16) -----------------
foreach ($idx in 14, 15, 16) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{C5A9637D-327B-4B98-B6D2-77B0DE613E2E}")
        }
    }
}

# --- 2. Swap the presentation's theme colors to the "Office" palette --
# (the deck's only reachable/active theme - driving the slide master,
# slides, notes master and handout master - currently carries the
# "Red Violet" / Integral palette; the edit swaps it for the stock
# "Office" palette that previously sat unused on the Notes Master theme)
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Index -> (scheme slot, RGB as 0xBBGGRR, i.e. R + G*256 + B*65536)
$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
